$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Task 12: status moved from "Open" to "In Progress"
$ws.Range("C12").Value = "In Progress"

# Task 13: status moved from "In Progress" to "Done" (and gets the
# green "Done" font colour used by the other completed rows)
$ws.Range("C13").Value = "Done"
$ws.Range("C13").Font.Color = $ws.Range("C4").Font.Color

# Cursor moved to C20
$ws.Range("C20").Select()
